# Weekly fruit/vegetable price update: insert a new week's worth of data
# (2 rows) at the top of the "Camote" price history, pushing all existing
# rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 62 (existing rows 62..128 shift down to 64..130)
$ws.Rows.Item(62).Insert()
$ws.Rows.Item(62).Insert()

# --- New row 62: Primera, $/caja 18 kilos ---
$ws.Cells.Item(62, 1).Value = 9
$ws.Cells.Item(62, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(62, 3).Value = 'Metropolitana'
$ws.Cells.Item(62, 4).Value = 44942
$ws.Cells.Item(62, 5).Value = 13
$ws.Cells.Item(62, 6).Value = 100114002
$ws.Cells.Item(62, 7).Value = 'Camote'
$ws.Cells.Item(62, 8).Value = 'Sin especificar'
$ws.Cells.Item(62, 9).Value = 'Primera'
$ws.Cells.Item(62, 10).Value = 700
$ws.Cells.Item(62, 11).Value = 17000
$ws.Cells.Item(62, 12).Value = 18000
$ws.Cells.Item(62, 13).Value = 17500
$ws.Cells.Item(62, 14).Value = '$/caja 18 kilos'
$ws.Cells.Item(62, 15).Value = 'Perú'
$ws.Cells.Item(62, 16).Value = 972
$ws.Cells.Item(62, 17).Value = 18
$ws.Cells.Item(62, 18).Value = 'Hortaliza'

# --- New row 63: Primera, $/malla 18 kilos ---
$ws.Cells.Item(63, 1).Value = 9
$ws.Cells.Item(63, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(63, 3).Value = 'Metropolitana'
$ws.Cells.Item(63, 4).Value = 44942
$ws.Cells.Item(63, 5).Value = 13
$ws.Cells.Item(63, 6).Value = 100114002
$ws.Cells.Item(63, 7).Value = 'Camote'
$ws.Cells.Item(63, 8).Value = 'Sin especificar'
$ws.Cells.Item(63, 9).Value = 'Primera'
$ws.Cells.Item(63, 10).Value = 610
$ws.Cells.Item(63, 11).Value = 16000
$ws.Cells.Item(63, 12).Value = 17000
$ws.Cells.Item(63, 13).Value = 16500
$ws.Cells.Item(63, 14).Value = '$/malla 18 kilos'
$ws.Cells.Item(63, 15).Value = 'Perú'
$ws.Cells.Item(63, 16).Value = 917
$ws.Cells.Item(63, 17).Value = 18
$ws.Cells.Item(63, 18).Value = 'Hortaliza'
